$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column A slightly (16.42578125 -> 15.42578125 in OOXML character-width
# units). The engine's ColumnWidth setter snaps to the nearest 1/6 character,
# same as real Excel's COM automation rounding, so 14.6 is the input that
# lands on the closest achievable stored width (15.5).
$ws.Columns.Item(1).ColumnWidth = 14.6

# Update the data column with the new values.
$ws.Range("A1").Value = 0.14727665526569922
$ws.Range("A2").Value = -0.0059999999697488704
$ws.Range("A3").Value = -0.0039999999741784364
$ws.Range("A4").Value = -0.0079999999519380083
$ws.Range("A5").Value = -0.0029999999737508887
$ws.Range("A6").Value = -0.001999999972168709
$ws.Range("A7").Value = -0.0099999999326501054
$ws.Range("A8").Value = -0.0099999999325488531
$ws.Range("A9").Value = -0.0019999999724999995
$ws.Range("A10").Value = -0.0019999999734157115
$ws.Range("A11").Value = -0.0029999999685532686
$ws.Range("A12").Value = -0.0025921005812654663
$ws.Range("A13").Value = 0.023470960514993422
$ws.Range("A14").Value = -0.0079999999465583116
$ws.Range("A15").Value = -0.00099999998098265763
$ws.Range("A16").Value = 0.010111267930441592
$ws.Range("A17").Value = -0.0019999999755633269
$ws.Range("A18").Value = -0.0039999999655497831
$ws.Range("A19").Value = -0.0039999999795266028
$ws.Range("A20").Value = -0.0039999999799178454
$ws.Range("A21").Value = 0.014048881039482097
$ws.Range("A22").Value = -0.0039999999799729125
$ws.Range("A23").Value = -0.00499999996763556
$ws.Range("A24").Value = -0.019999999890956133
$ws.Range("A25").Value = -0.019999999889500408
$ws.Range("A26").Value = -0.0024999999687622676
$ws.Range("A27").Value = -0.0024999999674095719
$ws.Range("A28").Value = -0.0019999999642124067
$ws.Range("A29").Value = -0.0069999999353198561
$ws.Range("A30").Value = -0.059999999673202353
$ws.Range("A31").Value = 0.023966696063721926
$ws.Range("A32").Value = -0.0099999999189339661
$ws.Range("A33").Value = -0.0039999999482684956
